$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = -21.846
$ws.Range("A12").Value = -21.696
$ws.Range("A18").Value = -22.23
$ws.Range("A37").Value = -20.217
$ws.Range("A55").Value = -22.164
$ws.Range("A68").Value = -21.696
$ws.Range("A77").Value = -20.505
$ws.Range("A78").Value = -19.951
